$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.022.16"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "3.565.44"
$ws.Range("E3").Value = "  -1.55%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'592.66"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").Value = "'197.37"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  -2.23%  "

$ws.Range("D8").Value = "'0.999"

$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("D10").Value = "'0.629"
$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("D11").Value = "'53.29"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("E12").Value = "  -4.46%  "

$ws.Range("D13").Value = "'9.31"
$ws.Range("E13").Value = "  -2.77%  "

$ws.Range("D14").Value = "4.126.41"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").Value = "'657.75"
$ws.Range("E15").Value = "  +9.59%  "

$ws.Range("D16").Value = "69.754.07"
$ws.Range("E16").Value = "  -1.27%  "

$ws.Range("D17").Value = "'12.63"
$ws.Range("E17").Value = "  -2.90%  "

$ws.Range("D18").Value = "3.560.62"
$ws.Range("E18").Value = "  -1.67%  "

$ws.Range("E19").Value = "  -0.94%  "

$ws.Range("D20").Value = "'18.53"
$ws.Range("E20").Value = "  -2.88%  "

$ws.Range("D21").Value = "'0.968"
$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("D22").Value = "'18.42"
$ws.Range("E22").Value = "  +3.22%  "

$ws.Range("D23").Value = "'5.36"
$ws.Range("E23").Value = "  +3.35%  "

$ws.Range("D24").Value = "'104.70"
$ws.Range("E24").Value = "  +2.82%  "

$ws.Range("D25").Value = "'4.42"
$ws.Range("E25").Value = "  -4.55%  "

$ws.Range("D26").Value = "'2.95"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("D27").Value = "'10.31"
$ws.Range("E27").Value = "  -4.30%  "

$ws.Range("D28").Value = "'9.67"
$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("D29").Value = "'33.52"
$ws.Range("E29").Value = "  -1.15%  "

$ws.Range("D30").Value = "'4.38"
$ws.Range("E30").Value = "  -6.57%  "

$ws.Range("D31").Value = "'6.85"
$ws.Range("E31").Value = "  -5.44%  "

$ws.Range("D32").Value = "'11.82"
$ws.Range("E32").Value = "  -3.87%  "

$ws.Range("D33").Value = "'0.112"
$ws.Range("E33").Value = "  -5.12%  "

$ws.Range("D34").Value = "'61.90"
$ws.Range("E34").Value = "  -2.45%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.756.61"
$ws.Range("E35").Value = "  -4.04%  "

$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").Value = "'3.75"
$ws.Range("E36").Value = "  +5.68%  "

$ws.Range("D37").Value = "0.0₃0820"
$ws.Range("E37").Value = "  -7.67%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").Value = "'515.96"
$ws.Range("E39").Value = "  -4.43%  "

$ws.Range("D40").Value = "'2.97"
$ws.Range("E40").Value = "  -4.83%  "

$ws.Range("E41").Value = "  -4.05%  "

$ws.Range("D42").Value = "'0.136"
$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").Value = "'35.09"
$ws.Range("E43").Value = "  -5.33%  "

$ws.Range("D44").Value = "'0.0453"
$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.41"
$ws.Range("E45").Value = "  +0.14%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.89"
$ws.Range("E46").Value = "  +0.82%  "

$ws.Range("E47").Value = "  -2.19%  "

$ws.Range("D48").Value = "'1.00"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("E49").Value = "  -2.59%  "

$ws.Range("B50").Value = "Jupiter"
$ws.Range("C50").Value = "https://coinranking.com/coin/qMgTxtv34+jupiter-jup"
$ws.Range("D50").Value = "'1.77"
$ws.Range("E50").Value = "  +19.22%  "

$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.72"
$ws.Range("E51").Value = "  +60.93%  "
